$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.981.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "'1.634.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("D5").Value = "'212.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").Value = "'23.48"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "'1.864.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "'1.630.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D15").Value = "'0.564"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "'65.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "'27.973.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'232.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "'7.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").Value = "'0.997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").Value = "'10.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.23%  "
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("E24").Value = "  -3.69%  "
$ws.Range("D25").Value = "'154.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("D27").Value = "'15.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "'1.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").Value = "'3.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.24%  "
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").Value = "'1.412.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'1.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.45%  "
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").Value = "'66.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("D47").Value = "'1.774.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").Value = "'88.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("E49").Value = "  -3.68%  "
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("E51").Value = "  -0.43%  "
